$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1985.1666
$ws.Range("I6").Value = 2282.2
$ws.Range("K6").Value = 6846.599999999999
$ws.Range("M6").Value = -6734.599999999999
$ws.Range("H15").Value = 93004.71000000001
$ws.Range("I15").Value = 93004.71000000001
$ws.Range("K15").Value = 279014.13
$ws.Range("M15").Value = -278845.13
$ws.Range("H39").Value = 1011.5833
$ws.Range("I39").Value = 1080.8182
$ws.Range("J39").Value = 250
$ws.Range("K39").Value = 3242.4546
$ws.Range("L39").Value = 750
$ws.Range("M39").Value = -2946.4546
$ws.Range("N39").Value = -1342
$ws.Range("H41").Value = 5291290
$ws.Range("I41").Value = 7936767.5
$ws.Range("J41").Value = 334.85715
$ws.Range("K41").Value = 7936767.5
$ws.Range("L41").Value = 334.85715
$ws.Range("M41").Value = -7936327.5
$ws.Range("N41").Value = -1214.85715
$ws.Range("H132").Value = 273397.22
$ws.Range("I132").Value = 358149.03
$ws.Range("J132").Value = 33267.082
$ws.Range("K132").Value = 1074447.09
$ws.Range("L132").Value = 99801.24600000001
$ws.Range("M132").Value = -1071917.09
$ws.Range("N132").Value = -104861.246
$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -50120
$ws.Range("H134").Value = 57000
$ws.Range("J134").Value = 57000
$ws.Range("L134").Value = 57000
$ws.Range("N134").Value = -67140
$ws.Range("H135").Value = 1195.0312
$ws.Range("I135").Value = 1074.7
$ws.Range("K135").Value = 9672.300000000001
$ws.Range("M135").Value = -7137.300000000001
$ws.Range("H137").Value = 1242.2222
$ws.Range("I137").Value = 694.7368
$ws.Range("J137").Value = 1854.1177
$ws.Range("K137").Value = 2084.2104
$ws.Range("L137").Value = 5562.3531
$ws.Range("M137").Value = 465.7896000000001
$ws.Range("N137").Value = -10662.3531

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2183.963
$ws.Range("I61").Value = 1125.55
$ws.Range("J61").Value = 5208
$ws.Range("K61").Value = 1125.55
$ws.Range("L61").Value = 5208
$ws.Range("M61").Value = -913.55
$ws.Range("N61").Value = -5632
$ws.Range("H74").Value = 1237.5588
$ws.Range("I74").Value = 1201.125
$ws.Range("K74").Value = 1201.125
$ws.Range("M74").Value = -327.125
$ws.Range("H77").Value = 1237.5588
$ws.Range("I77").Value = 1201.125
$ws.Range("K77").Value = 6005.625
$ws.Range("M77").Value = -1637.625
$ws.Range("H135").Value = 35000
$ws.Range("J135").Value = 35000
$ws.Range("L135").Value = 35000
$ws.Range("N135").Value = -45140
$ws.Range("H136").Value = 2183.963
$ws.Range("I136").Value = 1125.55
$ws.Range("J136").Value = 5208
$ws.Range("K136").Value = 3376.65
$ws.Range("L136").Value = 15624
$ws.Range("M136").Value = -826.6499999999996
$ws.Range("N136").Value = -20724
$ws.Range("H139").Value = 66600.28999999999
$ws.Range("J139").Value = 66600.28999999999
$ws.Range("L139").Value = 66600.28999999999
$ws.Range("N139").Value = -76880.28999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 42500
$ws.Range("J138").Value = 42500
$ws.Range("L138").Value = 42500
$ws.Range("N138").Value = -52780

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6251421
$ws.Range("I99").Value = 10417768
$ws.Range("K99").Value = 10417768
$ws.Range("M99").Value = -10416270
$ws.Range("H105").Value = 855.2
$ws.Range("I105").Value = 725.3333
$ws.Range("J105").Value = 1050
$ws.Range("K105").Value = 725.3333
$ws.Range("L105").Value = 1050
$ws.Range("M105").Value = 1021.6667
$ws.Range("N105").Value = -4544
$ws.Range("H126").Value = 6251421
$ws.Range("I126").Value = 10417768
$ws.Range("K126").Value = 31253304
$ws.Range("M126").Value = -31250834

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 173
$ws.Range("I7").Value = 173
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 519
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -407
$ws.Range("N7").ClearContents()
$ws.Range("H68").Value = 998.49
$ws.Range("I68").Value = 752.1429000000001
$ws.Range("J68").Value = 1417.9459
$ws.Range("K68").Value = 2256.4287
$ws.Range("L68").Value = 4253.8377
$ws.Range("M68").Value = -1445.4287
$ws.Range("N68").Value = -5875.8377
$ws.Range("H71").Value = 998.49
$ws.Range("I71").Value = 752.1429000000001
$ws.Range("J71").Value = 1417.9459
$ws.Range("K71").Value = 6769.2861
$ws.Range("L71").Value = 12761.5131
$ws.Range("M71").Value = -2713.2861
$ws.Range("N71").Value = -20873.5131
$ws.Range("H113").Value = 71429090
$ws.Range("I113").Value = 590
$ws.Range("J113").Value = 100000490
$ws.Range("K113").Value = 1770
$ws.Range("L113").Value = 300001470
$ws.Range("M113").Value = 400
$ws.Range("N113").Value = -300005810
$ws.Range("H131").Value = 2996.1428
$ws.Range("J131").Value = 3136.5593
$ws.Range("L131").Value = 9409.677899999999
$ws.Range("N131").Value = -19489.6779

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2687.5
$ws.Range("H83").Value = 2687.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3047.5
$ws.Range("I7").Value = 1712.5
$ws.Range("J7").Value = 3428.9285
$ws.Range("K7").Value = 1712.5
$ws.Range("L7").Value = 3428.9285
$ws.Range("M7").Value = -1600.5
$ws.Range("N7").Value = -3652.9285
$ws.Range("H126").Value = 3047.5
$ws.Range("I126").Value = 1712.5
$ws.Range("J126").Value = 3428.9285
$ws.Range("K126").Value = 5137.5
$ws.Range("L126").Value = 10286.7855
$ws.Range("M126").Value = -2667.5
$ws.Range("N126").Value = -15226.7855
$ws.Range("H136").Value = 3609.3235
$ws.Range("I136").Value = 1402.7778
$ws.Range("J136").Value = 12120.286
$ws.Range("K136").Value = 4208.3334
$ws.Range("L136").Value = 36360.858
$ws.Range("M136").Value = -1658.3334
$ws.Range("N136").Value = -41460.858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 383.5
$ws.Range("I100").Value = 275.25
$ws.Range("J100").Value = 600
$ws.Range("K100").Value = 550.5
$ws.Range("L100").Value = 1200
$ws.Range("M100").Value = -9.5
$ws.Range("N100").Value = -2282
$ws.Range("H132").Value = 11629615
$ws.Range("I132").Value = 16668012
$ws.Range("J132").Value = 2543.8462
$ws.Range("K132").Value = 50004036
$ws.Range("L132").Value = 7631.5386
$ws.Range("M132").Value = -50001506
$ws.Range("N132").Value = -12691.5386
$ws.Range("H136").Value = 8548097
$ws.Range("I136").Value = 11905330
$ws.Range("J136").Value = 2412.182
$ws.Range("K136").Value = 35715990
$ws.Range("L136").Value = 7236.545999999999
$ws.Range("M136").Value = -35713440
$ws.Range("N136").Value = -12336.546
